$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "CasesTab" row's Neo4j query (cell B2) referenced a `co:cohort` match
# that was only ever OPTIONAL-matched for the purpose of filtering, and the
# trailing `Cohort` output column was erroring out downstream. Drop that
# last RETURN line so the query only returns the fields that are actually
# backed by the MATCH/OPTIONAL MATCH clauses above it.
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Vizsla']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Reflect the author's resulting selection: the active cell moved to B2
# (previously B4 was selected with the grid scrolled so row 4 was the top
# visible row).
$ws.Range("B2").Select()
